# Auto-generated edit script for 北京-漫展信息.xlsx
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")

# --- simple counter bumps (rows unaffected by the row-29 insert) ---
$ws1.Range("F5").Value = 295
$ws1.Range("F6").Value = 427
$ws1.Range("F8").Value = 1932
$ws1.Range("F15").Value = 44
$ws1.Range("F18").Value = 4
$ws1.Range("F20").Value = 423
$ws1.Range("F24").Value = 6900
$ws1.Range("F25").Value = 7453
$ws1.Range("F26").Value = 29
$ws1.Range("F27").Value = 172

# --- insert a new row at position 29 (shifts old rows 29-48 down to 30-49) ---
$ws1.Rows.Item(29).Insert()
# copy column-A number formatting (border/bold/center) down onto the new row's A cell
$ws1.Range("A28").Copy()
$ws1.Range("A29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- populate the new row 29 ---
$ws1.Range("A29").Value = 28
$ws1.Range("B29").Value = "2024-10-02"
$ws1.Range("C29").Value = "北京·人气声优 内田秀 专场活动"
$ws1.Range("D29").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws1.Range("E29").Value = "2024.10.02 13:55-10.02 17:10"
$ws1.Range("F29").Value = 35
$ws1.Range("G29").Value = 458
$ws1.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=91678"
$ws1.Range("I29").Value = "//i0.hdslb.com/bfs/openplatform/202409/0aUkHD511725260741169.png"

# --- counter bumps on rows that were shifted down by the insert (new row numbers) ---
$ws1.Range("F31").Value = 233
$ws1.Range("F34").Value = 37
$ws1.Range("F36").Value = 1358
$ws1.Range("F37").Value = 7
$ws1.Range("F39").Value = 272
$ws1.Range("F40").Value = 661
$ws1.Range("F43").Value = 300
$ws1.Range("F47").Value = 100
$ws1.Range("F48").Value = 124

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value = 4

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 2530
$ws3.Range("F4").Value = 247
$ws3.Range("F5").Value = 103

# --- append new row 6 ---
# copy column-A number formatting (border/bold/center) onto the new row's A cell
$ws3.Range("A5").Copy()
$ws3.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws3.Range("A6").Value = 5
$ws3.Range("B6").Value = "2024-09-15"
$ws3.Range("C6").Value = "北京·夜境市集"
$ws3.Range("D6").Value = "南三环路 耕海大厦"
$ws3.Range("E6").Value = "2024.09.15 14:00-09.15 22:00"
$ws3.Range("F6").Value = 1
$ws3.Range("G6").Value = 128
$ws3.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=91663"
$ws3.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202408/ai0XaH8F1725011669001.png"

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 247
$ws4.Range("F7").Value = 103
$ws4.Range("F9").Value = 295
$ws4.Range("F11").Value = 427
$ws4.Range("F13").Value = 1932
$ws4.Range("F19").Value = 44
$ws4.Range("F21").Value = 423
$ws4.Range("F25").Value = 6900
$ws4.Range("F26").Value = 7453
$ws4.Range("F27").Value = 29
$ws4.Range("F28").Value = 172
$ws4.Range("F29").Value = 233
$ws4.Range("F30").Value = 37
$ws4.Range("F31").Value = 1358
$ws4.Range("F33").Value = 4
$ws4.Range("F34").Value = 272
$ws4.Range("F37").Value = 661
$ws4.Range("F43").Value = 300
$ws4.Range("F47").Value = 100

Write-Host "edit complete"
